$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''29.173.53'
$ws.Range("E2").Value = '  +0.29%  '

# Row 3
$ws.Range("D3").Value = '''1.829.93'
$ws.Range("E3").Value = '  -0.30%  '

# Row 4
$ws.Range("D4").Value = '''0.9985'

# Row 5
$ws.Range("D5").Value = '''242.81'
$ws.Range("E5").Value = '  -0.11%  '

# Row 6
$ws.Range("D6").Value = '''0.6207'
$ws.Range("E6").Value = '  +0.85%  '

# Row 7
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$ws.Range("D8").Value = '''0.07356'
$ws.Range("E8").Value = '  -1.45%  '

# Row 9
$ws.Range("D9").Value = '''0.2913'
$ws.Range("E9").Value = '  -0.17%  '

# Row 10
$ws.Range("D10").Value = '''23.21'
$ws.Range("E10").Value = '  +0.37%  '

# Row 11
$ws.Range("D11").Value = '''0.07660'
$ws.Range("E11").Value = '  -0.36%  '

# Row 12
$ws.Range("D12").Value = '''1.839.24'
$ws.Range("E12").Value = '  +0.10%  '

# Row 13
$ws.Range("D13").Value = '''4.971'
$ws.Range("E13").Value = '  -0.70%  '

# Row 14
$ws.Range("D14").Value = '''0.6688'
$ws.Range("E14").Value = '  -0.47%  '

# Row 15
$ws.Range("D15").Value = '''82.48'
$ws.Range("E15").Value = '  -0.21%  '

# Row 16
$ws.Range("D16").Value = '''0.000009015'
$ws.Range("E16").Value = '  -1.43%  '

# Row 17
$ws.Range("D17").Value = '''5.842'
$ws.Range("E17").Value = '  -1.43%  '

# Row 18
$ws.Range("D18").Value = '''29.169.19'
$ws.Range("E18").Value = '  +0.32%  '

# Row 19
$ws.Range("D19").Value = '''2.087.04'
$ws.Range("E19").Value = '  +0.43%  '

# Row 20
$ws.Range("D20").Value = '''235.35'
$ws.Range("E20").Value = '  +1.57%  '

# Row 21
$ws.Range("D21").Value = '''12.48'
$ws.Range("E21").Value = '  -1.38%  '

# Row 22
$ws.Range("D22").Value = '''0.9999'
$ws.Range("E22").Value = '  -0.22%  '

# Row 23
$ws.Range("D23").Value = '''7.361'
$ws.Range("E23").Value = '  +2.32%  '

# Row 24
$ws.Range("D24").Value = '''0.9992'
$ws.Range("E24").Value = '  -0.25%  '

# Row 25
$ws.Range("D25").Value = '''158.13'
$ws.Range("E25").Value = '  -0.96%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''8.549'
$ws.Range("E26").Value = '  +0.60%  '

# Row 27
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = '''0.1392'
$ws.Range("E27").Value = '  +0.26%  '

# Row 28
$ws.Range("D28").Value = '''17.61'
$ws.Range("E28").Value = '  -1.13%  '

# Row 29
$ws.Range("D29").Value = '''1.488'
$ws.Range("E29").Value = '  -0.46%  '

# Row 30
$ws.Range("D30").Value = '''0.05850'
$ws.Range("E30").Value = '  +5.81%  '

# Row 31
$ws.Range("D31").Value = '''4.090'
$ws.Range("E31").Value = '  -1.18%  '

# Row 32
$ws.Range("D32").Value = '''4.084'
$ws.Range("E32").Value = '  -1.76%  '

# Row 33
$ws.Range("E33").Value = '  +0.16%  '

# Row 34
$ws.Range("D34").Value = '''1.875'
$ws.Range("E34").Value = '  +2.12%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.140'
$ws.Range("E35").Value = '  +0.00%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.7263'
$ws.Range("E36").Value = '  -2.44%  '

# Row 37
$ws.Range("D37").Value = '''2.609'
$ws.Range("E37").Value = '  -1.99%  '

# Row 38
$ws.Range("D38").Value = '''2.853'
$ws.Range("E38").Value = '  +2.84%  '

# Row 39
$ws.Range("D39").Value = '''1.220.41'
$ws.Range("E39").Value = '  +0.78%  '

# Row 40
$ws.Range("D40").Value = '''0.01754'
$ws.Range("E40").Value = '  -1.52%  '

# Row 41
$ws.Range("D41").Value = '''6.236'
$ws.Range("E41").Value = '  -3.69%  '

# Row 42
$ws.Range("D42").Value = '''0.9085'
$ws.Range("E42").Value = '  +1.63%  '

# Row 43
$ws.Range("D43").Value = '''1.000'
$ws.Range("E43").Value = '  -0.07%  '

# Row 44
$ws.Range("D44").Value = '''1.992.47'
$ws.Range("E44").Value = '  +0.47%  '

# Row 45
$ws.Range("D45").Value = '''101.64'
$ws.Range("E45").Value = '  -0.42%  '

# Row 46
$ws.Range("D46").Value = '''65.43'
$ws.Range("E46").Value = '  -0.34%  '

# Row 47
$ws.Range("D47").Value = '''0.5041'
$ws.Range("E47").Value = '  -1.00%  '

# Row 48
$ws.Range("E48").Value = '  -3.03%  '

# Row 49
$ws.Range("D49").Value = '''9.161'
$ws.Range("E49").Value = '  +0.58%  '

# Row 50
$ws.Range("D50").Value = '''0.4027'
$ws.Range("E50").Value = '  -0.99%  '

# Row 51
$ws.Range("D51").Value = '''0.1129'
$ws.Range("E51").Value = '  +2.54%  '
